$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: fill in the previously-blank clock-out time and duration ---
# C22/D22 already carry the correct direct formatting (same as the rest of
# the data rows), so a plain value assignment is enough - no style changes
# are needed here.
$ws.Range("C22").Value = "12:59:28"
$ws.Range("D22").Value = "0.48 Hours"

# --- Row 23: brand-new timesheet entry ---
# These cells currently have no explicit formatting, so give them the same
# direct formatting as row 22 (PasteFormats keeps existing values/types
# untouched and reuses the same style index instead of creating a new one).
$ws.Range("A22:D22").Copy()
$ws.Range("A23:D23").PasteSpecial(-4122)  # xlPasteFormats

# Column A holds a literal "yyyy-mm-dd" text date. A plain .Value assignment
# would get silently auto-converted into a real date serial number by
# Excel, so instead enter it as a formula that evaluates to the literal
# string, then convert that formula to a static value in place. This keeps
# the cell as plain text without having to touch NumberFormat (which would
# otherwise register a brand-new, unused cell style).
$ws.Range("A23").Formula = "=""2026-02-07"""
$ws.Range("A23").Copy()
$ws.Range("A23").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B23").Value = "13:12:53"
$ws.Range("C23").Value = "13:40:48"
$ws.Range("D23").Value = "0.47 Hours"

# --- Row 24: totals row ---
$ws.Range("C22:D22").Copy()
$ws.Range("C24:D24").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C24").Value = "Total Duration:"
$ws.Range("D24").Value = "29 Hours"

$excel.CutCopyMode = 0
